# Updates the cryptocurrency price/volume table (and a few swapped rows)
# on Sheet1 to match the latest scrape, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "24.594.66"
$ws.Range("E2").Value = "  +3.16%  "

# Row 3
$ws.Range("D3").Value = "1.695.54"
$ws.Range("E3").Value = "  +1.76%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.72"
$ws.Range("E5").Value = "  +1.78%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  +0.04%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3935"
$ws.Range("E7").Value = "  +1.27%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4003"
$ws.Range("E8").Value = "  +0.79%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.532"
$ws.Range("E9").Value = "  +4.26%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9994"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.45"
$ws.Range("E11").Value = "  +4.13%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08748"
$ws.Range("E12").Value = "  +0.69%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.205"
$ws.Range("E13").Value = "  +7.05%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.20"
$ws.Range("E14").Value = "  +2.29%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.041"
$ws.Range("E15").Value = "  +9.82%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001316"
$ws.Range("E16").Value = "  +0.34%  "

# Row 17
$ws.Range("D17").Value = "1.694.92"
$ws.Range("E17").Value = "  +1.88%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "99.60"
$ws.Range("E18").Value = "  -0.08%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07072"
$ws.Range("E19").Value = "  +3.07%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.64"
$ws.Range("E20").Value = "  +2.58%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.942"
$ws.Range("E21").Value = "  +4.40%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9995"
$ws.Range("E22").Value = "  +0.15%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.11"
$ws.Range("E23").Value = "  +1.78%  "

# Row 24
$ws.Range("D24").Value = "24.583.18"
$ws.Range("E24").Value = "  +3.11%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.148"
$ws.Range("E25").Value = "  +10.01%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.351"
$ws.Range("E26").Value = "  +1.71%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.30"
$ws.Range("E27").Value = "  +2.40%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.66"
$ws.Range("E28").Value = "  +0.96%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.208"
$ws.Range("E29").Value = "  +1.03%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.98"
$ws.Range("E30").Value = "  +3.61%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.527"
$ws.Range("E31").Value = "  +10.97%  "

# Row 32
$ws.Range("D32").Value = "1.880.37"
$ws.Range("E32").Value = "  +1.61%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.082"
$ws.Range("E33").Value = "  -3.71%  "

# Row 34
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08559"
$ws.Range("E34").Value = "  +0.61%  "

# Row 35
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.241"
$ws.Range("E35").Value = "  +9.12%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.34"
$ws.Range("E36").Value = "  +8.80%  "

# Row 37
$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2722"
$ws.Range("E37").Value = "  +2.71%  "

# Row 38
$ws.Range("B38").Value = "WEMIXTOKEN"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.930"
$ws.Range("E38").Value = "  -0.27%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.46"
$ws.Range("E39").Value = "  -0.39%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02742"
$ws.Range("E40").Value = "  +8.90%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09051"
$ws.Range("E41").Value = "  +2.90%  "

# Row 42
$ws.Range("E42").Value = "  +1.23%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7653"
$ws.Range("E43").Value = "  +1.00%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7173"
$ws.Range("E44").Value = "  +1.78%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.50"
$ws.Range("E45").Value = "  +2.35%  "

# Row 46
$ws.Range("E46").Value = "  +3.27%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.208"
$ws.Range("E47").Value = "  +2.53%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9993"
$ws.Range("E48").Value = "  +0.01%  "

# Row 49
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.02"
$ws.Range("E49").Value = "  +0.69%  "

# Row 50
$ws.Range("B50").Value = "Flow"
$ws.Range("C50").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.324"
$ws.Range("E50").Value = "  +8.82%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07985"
$ws.Range("E51").Value = "  +2.34%  "
